# Applies the "added culture data browsing functions" edit to dna_data.xlsx
# - Updates the shared Help-text banner on both sheets with a bold warning run.
# - Renames several Help-sheet "Database Field" / "Database Table" helper
#   strings (rows 11-28, columns B and C) to the new naming scheme.
# - Makes the DNA sheet the active tab/sheet (was Help), and normalizes both
#   sheets' selections to A2.

$wb = $excel.ActiveWorkbook
$wsHelp = $wb.Worksheets.Item("Help")
$wsDNA  = $wb.Worksheets.Item("DNA")

# ---------------------------------------------------------------------------
# 1. Shared "help banner" text in A2 on both sheets, with a bold, sz=12 run
#    appended as a strong caution after the original sentence. Both sheets
#    originally pointed at the very same shared string, so build the rich
#    text once on Help!A2 and copy/paste-values it onto DNA!A2 so the two
#    cells keep sharing a single (deduped) shared-string entry instead of
#    each getting their own private copy of the rich text.
# ---------------------------------------------------------------------------
$bannerPlain = "This file will be used to upload data to the NelsonDB. "
$bannerBold  = "FIELD NAMES CAN BE ADDED, BUT SHOULD BE DONE SPARINGLY"
$banner = $bannerPlain + $bannerBold
$boldStart = $bannerPlain.Length + 1

$wsHelp.Range("A2").Value = $banner
$wsHelp.Range("A2").Characters($boldStart, $bannerBold.Length).Font.Bold = $true
$wsHelp.Range("A2").Characters($boldStart, $bannerBold.Length).Font.Size = 12

$wsHelp.Range("A2").Copy()
$wsDNA.Range("A2").PasteSpecial(-4163)
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 2. Help sheet "DNA Info" rows (11-16): rename the table name in column B
#    from lab_obsother -> lab_obsdna, and drop the "dna_" prefix on the field
#    names in column C.
# ---------------------------------------------------------------------------
$wsHelp.Range("B11:B16").Value = "lab_obsdna"

$wsHelp.Range("C11").Value = "dna_id"
$wsHelp.Range("C12").Value = "extraction_method"
$wsHelp.Range("C13").Value = "date"
$wsHelp.Range("C14").Value = "tube_id"
$wsHelp.Range("C15").Value = "tube_type"
$wsHelp.Range("C16").Value = "comments"

# ---------------------------------------------------------------------------
# 3. Help sheet "Source X ID" rows (18-28): rename the table name in column B
#    from lab_obsother -> lab_obstracker, and rewrite the mapping formulas in
#    column C to the new obs_*_id = Obs*(...).id convention.
# ---------------------------------------------------------------------------
$wsHelp.Range("B18:B28").Value = "lab_obstracker"

$wsHelp.Range("C18").Value = "stock_id = Stock(seed_id).id"
$wsHelp.Range("C19").Value = "isolate_id = Isolate(isolate_id).id"
$wsHelp.Range("C20").Value = "obs_row_id = ObsRow(row_id).id"
$wsHelp.Range("C21").Value = "obs_plant_id = ObsPlant(plant_id).id"
$wsHelp.Range("C22").Value = "obs_well_id = ObsWellr(well_id).id"
$wsHelp.Range("C23").Value = "obs_microbe_id = ObsMicrobe(microbe_id).id"
$wsHelp.Range("C24").Value = "obs_culture_id = ObsCulture(culture_id).id"
$wsHelp.Range("C25").Value = "obs_tissue_id = ObsTissue(tissue_id).id"
$wsHelp.Range("C26").Value = "obs_sample_id = ObsSample(sample_id).id"
$wsHelp.Range("C27").Value = "obs-plate_id = ObsPlate(plate_id).id"
$wsHelp.Range("C28").Value = "obs_dna_id = ObsDNA(dna_id).id"

# ---------------------------------------------------------------------------
# 4. View state: DNA becomes the active/selected tab (was Help); both sheets'
#    selection cursor moves to A2 (was C7 on Help, C12 on DNA).
# ---------------------------------------------------------------------------
$wsHelp.Activate()
$wsHelp.Range("A2").Select()

$wsDNA.Activate()
$wsDNA.Range("A2").Select()
